# Commit: "adding sample id column in dqo accuracy check tables #15"
#
# For this workbook (inst/extdata/ExampleDQOAccuracy.xlsx, sheet "Accuracy"),
# the resulting content changes are:
#   - "Conductivity" parameter renamed to "Sp Conductance" (rows 6-7)
#   - Several previously blank ("na") / placeholder Field Blank & Lab Blank
#     cells (columns H & I) are filled in with their proper detection-limit
#     values (BDL or a numeric/text threshold), matching the values already
#     used by sibling rows for the same parameter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Parameter name correction: Conductivity -> Sp Conductance
$ws.Range("A6").Value = "Sp Conductance"
$ws.Range("A7").Value = "Sp Conductance"

# DO (rows 4-5): Field Blank / Lab Blank
$ws.Range("H5").Value = "< 0.1"
$ws.Range("I5").Value = "< 0.1"

# Sp Conductance (rows 6-7): Field Blank / Lab Blank
$ws.Range("H7").Value = "< 25"
$ws.Range("I7").Value = "< 25"

# TSS (rows 8-9): Field Blank / Lab Blank
$ws.Range("H9").Value = "BDL"
$ws.Range("I9").Value = "BDL"

# TP (rows 10-11): Field Blank / Lab Blank
$ws.Range("H11").Value = "BDL"
$ws.Range("I11").Value = "BDL"

# Ortho P (rows 12-13): Field Blank / Lab Blank
$ws.Range("H13").Value = "BDL"
$ws.Range("I13").Value = "BDL"

# Nitrate (row 14): Field Blank / Lab Blank
$ws.Range("H14").Value = "BDL"
$ws.Range("I14").Value = "BDL"

# Ammonia (row 15): Field Blank / Lab Blank
$ws.Range("H15").Value = "BDL"
$ws.Range("I15").Value = "BDL"

# Chl a (rows 17-18): Field Blank / Lab Blank
$ws.Range("H17").Value = "BDL"
$ws.Range("I17").Value = "BDL"
$ws.Range("H18").Value = "BDL"
$ws.Range("I18").Value = "BDL"
